$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D range to Text format so numeric-looking values
# (e.g. "1.005", "14.76") are stored as text strings, matching the
# original inlineStr cell type, not auto-converted to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.977.64"
$ws.Range("D3").Value = "1.825.39"
$ws.Range("D4").Value = "1.005"
$ws.Range("D5").Value = "311.67"
$ws.Range("D7").Value = "0.4622"
$ws.Range("D8").Value = "0.3704"
$ws.Range("D10").Value = "0.8742"
$ws.Range("D13").Value = "1.805.19"
$ws.Range("D14").Value = "5.337"
$ws.Range("D15").Value = "6.542"
$ws.Range("D16").Value = "91.19"
$ws.Range("D18").Value = "0.000008859"
$ws.Range("D19").Value = "1.005"
$ws.Range("D20").Value = "14.76"
$ws.Range("D21").Value = "26.997.48"
$ws.Range("D22").Value = "5.099"
$ws.Range("D23").Value = "10.52"
$ws.Range("D24").Value = "2.068.07"
$ws.Range("D25").Value = "152.79"
$ws.Range("D26").Value = "1.846"
$ws.Range("D27").Value = "18.40"
$ws.Range("D28").Value = "2.035"
$ws.Range("D29").Value = "5.132"
$ws.Range("D30").Value = "115.48"
$ws.Range("D32").Value = "2.964"
$ws.Range("D33").Value = "0.7268"
$ws.Range("D34").Value = "4.430"
$ws.Range("D36").Value = "2.468"
$ws.Range("D37").Value = "0.01947"
$ws.Range("D39").Value = "0.05224"
$ws.Range("D42").Value = "0.5147"
$ws.Range("D43").Value = "0.1619"
$ws.Range("D45").Value = "0.4839"
$ws.Range("D47").Value = "10.17"
$ws.Range("D48").Value = "102.81"
$ws.Range("D49").Value = "1.631"
$ws.Range("D50").Value = "0.06196"
$ws.Range("D51").Value = "64.62"

# Restore the default (Normal) style on column D so no extra
# per-cell style index is left behind, matching the original file.
$dRange.Style = "Normal"

# Column E values are unambiguous text (contain "%", spaces, sign)
# so they remain text without needing a NumberFormat override.
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  -3.98%  "
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("E51").Value = "  +0.18%  "
